$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated metrics (access/citation/altmetric counts) for paper rows
# Each entry: row number, new F (access), new G (citation) or "" if unchanged, new H (altmetric) or "" if unchanged
$updates = @(
    ,@(2, 95, "", "")
    ,@(3, 55, "", "")
    ,@(4, 87, "", "")
    ,@(5, 46, "", 4)
    ,@(6, 37, "", "")
    ,@(7, 106, "", "")
    ,@(8, 468, "", "")
    ,@(9, 135, "", "")
    ,@(10, 93, "", "")
    ,@(11, 98, "", "")
    ,@(12, 266, "", "")
    ,@(13, 221, "", "")
    ,@(14, 1373, "", "")
    ,@(15, 232, "", "")
    ,@(16, 233, "", 6)
    ,@(17, 591, "", "")
    ,@(18, 190, "", 8)
    ,@(19, 418, "", "")
    ,@(20, 802, "", "")
    ,@(21, 217, "", "")
    ,@(22, 924, "", "")
    ,@(23, 1039, "", "")
    ,@(24, 577, "", "")
    ,@(25, 324, "", "")
    ,@(26, 475, "", "")
    ,@(27, 1075, "", "")
    ,@(28, 462, "", "")
    ,@(29, 311, "", "")
    ,@(30, 308, "", "")
    ,@(31, 326, "", "")
    ,@(32, 568, "", "")
    ,@(33, 964, "", "")
    ,@(34, 560, "", "")
    ,@(35, 765, "", "")
    ,@(36, 624, "", "")
    ,@(37, 588, "", "")
    ,@(38, 506, "", 12)
    ,@(39, 1135, "", "")
    ,@(40, 595, "", "")
    ,@(41, 689, "", "")
    ,@(42, 2778, "", "")
    ,@(43, 1079, "", "")
    ,@(44, 1371, "", "")
    ,@(45, 565, "", "")
    ,@(46, 2554, 18, "")
    ,@(47, 742, "", "")
    ,@(48, 1003, "", "")
    ,@(49, 2838, "", "")
    ,@(50, 1408, "", "")
    ,@(51, 1040, "", "")
    ,@(52, 936, "", "")
    ,@(53, 703, "", "")
    ,@(54, 703, "", "")
    ,@(55, 1129, "", "")
    ,@(56, 2645, "", "")
    ,@(57, 765, "", "")
    ,@(58, 881, "", "")
    ,@(59, 721, "", "")
    ,@(60, 1907, "", "")
    ,@(61, 588, "", "")
    ,@(62, 802, "", "")
    ,@(63, 679, "", "")
    ,@(64, 913, "", "")
    ,@(65, 1133, "", "")
    ,@(66, 653, "", "")
    ,@(67, 969, "", "")
    ,@(68, 2235, "", "")
    ,@(69, 2939, "", "")
    ,@(70, 2394, "", "")
    ,@(71, 5672, "", "")
    ,@(72, 845, "", "")
    ,@(73, 3188, "", "")
    ,@(74, 941, "", "")
    ,@(75, 4853, "", "")
    ,@(76, 4469, "", "")
    ,@(77, 1564, "", "")
    ,@(78, 915, "", "")
    ,@(79, 851, "", "")
    ,@(80, 1889, "", "")
    ,@(81, 8483, "", "")
    ,@(82, 987, "", "")
    ,@(83, 1200, "", "")
    ,@(84, 1322, "", "")
    ,@(85, 4937, "", "")
    ,@(86, 1446, "", "")
    ,@(87, 1496, "", "")
    ,@(88, 1824, "", "")
    ,@(89, 1401, "", "")
    ,@(90, 1487, "", "")
    ,@(91, 1664, "", "")
    ,@(92, 1497, "", "")
    ,@(93, 3869, "", "")
    ,@(94, 4690, "", "")
    ,@(95, 1964, "", "")
    ,@(96, 1627, "", "")
    ,@(97, 1713, "", "")
    ,@(98, 4026, "", "")
    ,@(99, 1494, "", "")
    ,@(101, 2277, "", "")
    ,@(102, 2025, "", "")
    ,@(103, 1494, "", "")
    ,@(104, 6196, "", "")
    ,@(105, 1691, "", "")
    ,@(106, 1371, "", "")
    ,@(107, 3190, "", "")
    ,@(108, 4211, "", "")
    ,@(109, 1776, "", "")
    ,@(110, 3655, "", "")
    ,@(111, 2128, "", "")
    ,@(112, 2112, "", "")
    ,@(113, 2362, "", "")
    ,@(114, 3091, "", "")
    ,@(115, 1725, "", "")
    ,@(116, 2644, "", "")
    ,@(117, 2069, "", "")
    ,@(118, 1644, "", "")
    ,@(119, 2085, "", "")
    ,@(120, 4719, "", "")
    ,@(121, 1423, "", "")
    ,@(122, 2601, "", "")
    ,@(123, 1673, "", "")
    ,@(124, 6988, 84, "")
    ,@(125, 1519, "", "")
    ,@(126, 1397, "", "")
    ,@(127, 1415, "", "")
    ,@(128, 2370, "", "")
    ,@(129, 2088, "", "")
    ,@(130, 1283, "", "")
    ,@(131, 3723, "", "")
    ,@(132, 1985, "", "")
    ,@(133, 3560, "", "")
    ,@(134, 2584, "", 6)
    ,@(135, 1565, "", "")
    ,@(136, 3671, 32, "")
    ,@(137, 2865, "", "")
    ,@(138, 2569, "", "")
    ,@(139, 1720, "", "")
    ,@(140, 3222, "", "")
    ,@(141, 4858, "", "")
    ,@(142, 5690, "", "")
    ,@(143, 1662, "", "")
    ,@(144, 2344, "", "")
    ,@(145, 2531, "", "")
    ,@(146, 2457, "", "")
    ,@(147, 2392, "", "")
    ,@(148, 2440, "", "")
    ,@(149, 2076, "", "")
    ,@(150, 5075, "", "")
    ,@(151, 2337, "", "")
    ,@(152, 2185, 15, "")
    ,@(153, 2671, "", "")
    ,@(154, 6014, "", "")
    ,@(155, 5426, "", "")
    ,@(156, 5375, "", "")
    ,@(157, 5660, "", "")
    ,@(158, 4749, "", "")
    ,@(159, 1774, "", "")
    ,@(160, 2477, "", "")
    ,@(161, 4132, "", "")
    ,@(162, 6125, "", "")
    ,@(163, 5031, "", "")
    ,@(164, 6274, "", "")
    ,@(165, 3110, "", "")
    ,@(166, 8064, 7, "")
    ,@(167, 6775, 17, "")
)

foreach ($u in $updates) {
    $row = $u[0]
    $fVal = $u[1]
    $gVal = $u[2]
    $hVal = $u[3]
    $ws.Cells.Item($row, 6).Value = $fVal
    if ($gVal -ne "") {
        $ws.Cells.Item($row, 7).Value = $gVal
    }
    if ($hVal -ne "") {
        $ws.Cells.Item($row, 8).Value = $hVal
    }
}
